$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 125000420
$ws.Range("I5").Value = 123.5
$ws.Range("J5").Value = 250000720
$ws.Range("K5").Value = 123.5
$ws.Range("L5").Value = 250000720
$ws.Range("M5").Value = -8.5
$ws.Range("N5").Value = -250000950

$ws.Range("H17").Value = 1870.5
$ws.Range("J17").Value = 2064.6
$ws.Range("L17").Value = 6193.799999999999
$ws.Range("N17").Value = -6529.799999999999

$ws.Range("H33").Value = 204.55556
$ws.Range("J33").Value = 391.25
$ws.Range("L33").Value = 391.25
$ws.Range("N33").Value = -849.25

$ws.Range("H55").Value = 100000180
$ws.Range("I55").Value = 142857340
$ws.Range("K55").Value = 142857340
$ws.Range("M55").Value = -142857126

$ws.Range("H70").Value = 15899.857
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 15899.857
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H87").Value = 26666.666

$ws.Range("H90").Value = 26666.666

$ws.Range("H103").Value = 895.6667
$ws.Range("I103").Value = 226.33333
$ws.Range("J103").Value = 1029.5333
$ws.Range("K103").Value = 678.99999
$ws.Range("L103").Value = 3088.5999
$ws.Range("M103").Value = -92.99999000000003
$ws.Range("N103").Value = -4260.5999

$ws.Range("H116").Value = 8097.364
$ws.Range("I116").Value = 8629.200000000001
$ws.Range("J116").Value = 2779
$ws.Range("K116").Value = 8629.200000000001
$ws.Range("L116").Value = 2779
$ws.Range("M116").Value = -5187.200000000001
$ws.Range("N116").Value = -9663

$ws.Range("H125").Value = 3544.4666
$ws.Range("I125").Value = 2333.6667
$ws.Range("J125").Value = 4351.6665
$ws.Range("K125").Value = 21003.0003
$ws.Range("L125").Value = 39164.9985
$ws.Range("M125").Value = -18543.0003
$ws.Range("N125").Value = -44084.9985

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 24888.666
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H61").Value = 1828.8889
$ws.Range("I61").Value = 1041.9231
$ws.Range("K61").Value = 1041.9231
$ws.Range("M61").Value = -829.9231

$ws.Range("H63").Value = 1065.5
$ws.Range("I63").Value = 463.33334
$ws.Range("J63").Value = 1968.75
$ws.Range("K63").Value = 463.33334
$ws.Range("L63").Value = 1968.75
$ws.Range("M63").Value = 222.66666
$ws.Range("N63").Value = -3340.75

$ws.Range("H66").Value = 1065.5
$ws.Range("I66").Value = 463.33334
$ws.Range("J66").Value = 1968.75
$ws.Range("K66").Value = 2316.6667
$ws.Range("L66").Value = 9843.75
$ws.Range("M66").Value = 1115.3333
$ws.Range("N66").Value = -16707.75

$ws.Range("H97").Value = 3259.3845
$ws.Range("I97").Value = 689.6111
$ws.Range("J97").Value = 9041.375
$ws.Range("K97").Value = 689.6111
$ws.Range("L97").Value = 9041.375
$ws.Range("M97").Value = -193.6111
$ws.Range("N97").Value = -10033.375

$ws.Range("H132").Value = 4445.3335
$ws.Range("I132").Value = 4246.7915
$ws.Range("J132").Value = 4974.778
$ws.Range("K132").Value = 12740.3745
$ws.Range("L132").Value = 14924.334
$ws.Range("M132").Value = -10210.3745
$ws.Range("N132").Value = -19984.334

$ws.Range("H136").Value = 1828.8889
$ws.Range("I136").Value = 1041.9231
$ws.Range("K136").Value = 3125.7693
$ws.Range("M136").Value = -575.7692999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15290.5
$ws.Range("I82").Value = 4348.7
$ws.Range("J82").Value = 69999.5
$ws.Range("K82").Value = 4348.7
$ws.Range("L82").Value = 69999.5
$ws.Range("M82").Value = -3965.7
$ws.Range("N82").Value = -70765.5

$ws.Range("H85").Value = 15290.5
$ws.Range("I85").Value = 4348.7
$ws.Range("J85").Value = 69999.5
$ws.Range("K85").Value = 4348.7
$ws.Range("L85").Value = 69999.5
$ws.Range("M85").Value = -3022.7
$ws.Range("N85").Value = -72651.5

$ws.Range("H86").Value = 3559.1
$ws.Range("I86").Value = 3098.8572
$ws.Range("J86").Value = 4633
$ws.Range("K86").Value = 3098.8572
$ws.Range("L86").Value = 4633
$ws.Range("M86").Value = -1975.8572
$ws.Range("N86").Value = -6879

$ws.Range("H89").Value = 3559.1
$ws.Range("I89").Value = 3098.8572
$ws.Range("J89").Value = 4633
$ws.Range("K89").Value = 15494.286
$ws.Range("L89").Value = 23165
$ws.Range("M89").Value = -9878.286
$ws.Range("N89").Value = -34397

$ws.Range("H94").Value = 5760.6313
$ws.Range("I94").Value = 520.7059
$ws.Range("K94").Value = 520.7059
$ws.Range("M94").Value = -69.70590000000004

$ws.Range("H107").Value = 1425.5834
$ws.Range("I107").Value = 1047
$ws.Range("J107").Value = 2561.3333
$ws.Range("K107").Value = 1047
$ws.Range("L107").Value = 2561.3333
$ws.Range("M107").Value = 873
$ws.Range("N107").Value = -6401.3333

$ws.Range("H134").Value = 7621.1064
$ws.Range("I134").Value = 7776.75
$ws.Range("K134").Value = 23330.25
$ws.Range("M134").Value = -20795.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3041.04
$ws.Range("I122").Value = 2548.647
$ws.Range("K122").Value = 7645.941
$ws.Range("M122").Value = -5195.941

$ws.Range("H132").Value = 2503.6738
$ws.Range("I132").Value = 1919.909
$ws.Range("K132").Value = 5759.727000000001
$ws.Range("M132").Value = -3229.727000000001

$ws.Range("H134").Value = 1414.1569
$ws.Range("I134").Value = 1313.919
$ws.Range("K134").Value = 3941.757000000001
$ws.Range("M134").Value = -1406.757000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 4000
$ws.Range("I19").Value = 4000
$ws.Range("K19").Value = 12000
$ws.Range("M19").Value = -11826

$ws.Range("H69").Value = 4921.4287
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 4992.3076
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 14976.9228
$ws.Range("M69").Value = -11189
$ws.Range("N69").Value = -16598.9228

$ws.Range("H72").Value = 4921.4287
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 4992.3076
$ws.Range("K72").Value = 36000
$ws.Range("L72").Value = 44930.7684
$ws.Range("M72").Value = -31944
$ws.Range("N72").Value = -53042.7684

$ws.Range("H132").Value = 3629.7097
$ws.Range("I132").Value = 2094.8
$ws.Range("J132").Value = 3924.8845
$ws.Range("K132").Value = 18853.2
$ws.Range("L132").Value = 35323.9605
$ws.Range("M132").Value = -16323.2
$ws.Range("N132").Value = -40383.9605

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 3509999.5
$ws.Range("I20").Value = 3509999.5
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 3509999.5
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -3509754.5
$ws.Range("N20").ClearContents()

$ws.Range("H70").Value = 11262.728
$ws.Range("J70").Value = 7597
$ws.Range("L70").Value = 7597
$ws.Range("N70").Value = -8137

$ws.Range("H73").Value = 11262.728
$ws.Range("J73").Value = 7597
$ws.Range("L73").Value = 7597
$ws.Range("N73").Value = -9469

$ws.Range("H102").Value = 7918.049
$ws.Range("I102").Value = 8175.3423
$ws.Range("K102").Value = 8175.3423
$ws.Range("M102").Value = -6553.3423

$ws.Range("H113").Value = 2067.6155
$ws.Range("I113").Value = 2079.2
$ws.Range("K113").Value = 2079.2
$ws.Range("M113").Value = 90.80000000000018

$ws.Range("H132").Value = 2813.8545
$ws.Range("I132").Value = 2272.2144
$ws.Range("J132").Value = 4563.769
$ws.Range("K132").Value = 6816.6432
$ws.Range("L132").Value = 13691.307
$ws.Range("M132").Value = -4286.6432
$ws.Range("N132").Value = -18751.307

$ws.Range("H136").Value = 46940.45
$ws.Range("J136").Value = 46940.45
$ws.Range("L136").Value = 140821.35
$ws.Range("N136").Value = -145921.35

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1852
$ws.Range("I22").Value = 1156.3334
$ws.Range("J22").Value = 2373.75
$ws.Range("K22").Value = 1156.3334
$ws.Range("L22").Value = 2373.75
$ws.Range("M22").Value = -861.3334
$ws.Range("N22").Value = -2963.75

$ws.Range("H27").Value = 1852
$ws.Range("I27").Value = 1156.3334
$ws.Range("J27").Value = 2373.75
$ws.Range("K27").Value = 1156.3334
$ws.Range("L27").Value = 2373.75
$ws.Range("M27").Value = -1049.3334
$ws.Range("N27").Value = -2587.75

$ws.Range("H40").Value = 2400.8462
$ws.Range("I40").Value = 2205.5789
$ws.Range("K40").Value = 2205.5789
$ws.Range("M40").Value = -2069.5789

$ws.Range("H46").Value = 1386.8889
$ws.Range("I46").Value = 1191.2
$ws.Range("J46").Value = 1631.5
$ws.Range("K46").Value = 1191.2
$ws.Range("L46").Value = 1631.5
$ws.Range("M46").Value = -1003.2
$ws.Range("N46").Value = -2007.5

$ws.Range("H100").Value = 3882.8333
$ws.Range("I100").Value = 3074.5
$ws.Range("J100").Value = 5499.5
$ws.Range("K100").Value = 3074.5
$ws.Range("L100").Value = 5499.5
$ws.Range("M100").Value = -2533.5
$ws.Range("N100").Value = -6581.5

$ws.Range("H122").Value = 2810.484
$ws.Range("I122").Value = 2312.1333
$ws.Range("J122").Value = 3277.6875
$ws.Range("K122").Value = 6936.3999
$ws.Range("L122").Value = 9833.0625
$ws.Range("M122").Value = -4486.3999
$ws.Range("N122").Value = -14733.0625

$ws.Range("H132").Value = 3728.762
$ws.Range("J132").Value = 4280.4614
$ws.Range("L132").Value = 12841.3842
$ws.Range("N132").Value = -17901.3842

$ws.Range("H136").Value = 3310.6875
$ws.Range("I136").Value = 3166.1333
$ws.Range("K136").Value = 9498.3999
$ws.Range("M136").Value = -6948.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11750.125
$ws.Range("J81").Value = 4000
$ws.Range("L81").Value = 8000
$ws.Range("N81").Value = -10122

$ws.Range("H84").Value = 11750.125
$ws.Range("J84").Value = 4000
$ws.Range("L84").Value = 40000
$ws.Range("N84").Value = -50608
